# Updates cryptos list values (price and 1h volume % change) for rows 2-51.
# Rows 31 and 32 also swap Coin name / Link (Filecoin <-> InternetComputer(DFINITY)).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.982.71"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.42%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'1.826.83"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.43%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'0.9974"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +0.61%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'241.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.47%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'0.6153"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -1.60%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'0.9994"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.44%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.07364"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.84%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.2932"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +0.69%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'22.93"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +0.15%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.07642"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.26%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'1.845.58"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +1.85%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'4.979"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.16%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'0.6686"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +0.73%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'82.39"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.24%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'0.000008996"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -6.21%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'5.872"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -2.24%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'29.023.40"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.63%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'2.085.67"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +1.73%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'234.75"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +5.12%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'12.65"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +1.13%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'0.9995"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.51%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'7.164"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +1.32%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'0.9982"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.49%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'158.46"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.19%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'0.1421"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +1.50%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'8.456"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +0.05%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'17.76"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -0.07%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'1.493"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +0.18%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'0.05558"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +2.27%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "'4.110"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +0.34%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Value = "'4.091"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +1.47%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'1.208"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +1.57%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'1.839"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -0.12%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'0.7401"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +0.13%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'1.133"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +0.34%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'2.642"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +1.76%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'2.774"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +1.40%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'0.01771"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.02%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'1.204.44"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -1.65%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'6.330"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -4.65%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'0.8962"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +0.45%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'0.9982"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +0.30%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'101.10"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.14%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'1.989.81"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +1.94%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'64.77"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.04%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'0.00000000122"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.41%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'0.5078"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +0.45%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'0.4039"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +0.50%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'9.071"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +1.66%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'0.05806"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +0.51%  "
$ws.Range("E51").Style = "Normal"
